# The commit deletes worksheet row 3 (which held the sample "Packs" /
# "dasjalsa" test data) entirely, shifting every row below it up by one.
# This is why the sheet's used dimension shrinks from A1:K96 to A1:K95,
# why the old row 4 (unit-price formulas) becomes the new row 3, why the
# old row 5 (SUM totals) becomes the new row 4, and why the very last
# row (96) disappears once everything has shifted up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire 3rd row, shifting rows 4:96 up to 3:95.
$ws.Rows.Item(3).Delete()
